{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\nfor (let i = 0; i < items.length; i++) {\n  const text = items[i].text;\n\n  if (text.indexOf(\"Make a one hot map of venue nearby and compare that with user's choice.\") !== -1) {\n    items[i].insertText(\n      \" Use eucleadian distance for finding the 4 or 5 best fit neighbourhoods.\",\n      \"End\"\n    );\n  } else if (\n    text.indexOf(\n      \"I would also be using the rent data from the following URL, to take into account affordibility of the neighbourhood.\"\n    ) !== -1\n  ) {\n    items[i].insertText(\n      \" Avg rent column from the link would be used to define the rent of the neighourhood\",\n      \"End\"\n    );\n  } else if (\n    text.indexOf(\"Priority should be given to proximity, facilities, affordibility in that order\") !== -1\n  ) {\n    items[i].insertText(\n      \"From the already shorlisted neighbourhoods, I'll check which one has the least amount of rent.\",\n      \"Replace\"\n    );\n    items[i].insertParagraph(\n      \"Priority should be given to proximity (will be taken care by foursquare radius variable), facilities (will be taken care by best fit neighourhoods), affordibility in that order (will be taken care by the rent data)\",\n      \"After\"\n    );\n  }\n}\n\nawait context.sync();\n", "ps1": "$wdCollapseEnd = 0\n\n$d = $word.ActiveDocument\n\n# --- Edit 1: append sentence to the \"one hot map\" paragraph ---\n$rng1 = $d.Content\n$rng1.Find.Text = \"Make a one hot map of venue nearby and compare that with user's choice.\"\n$rng1.Find.Execute() | Out-Null\n$rng1.Collapse($wdCollapseEnd)\n$rng1.InsertAfter(\" Use eucleadian distance for finding the 4 or 5 best fit neighbourhoods.\")\n\n# --- Edit 2: append sentence to the \"rent data\" paragraph ---\n$rng2 = $d.Content\n$rng2.Find.Text = \"I would also be using the rent data from the following URL, to take into account affordibility of the neighbourhood.\"\n$rng2.Find.Execute() | Out-Null\n$rng2.Collapse($wdCollapseEnd)\n$rng2.InsertAfter(\" Avg rent column from the link would be used to define the rent of the neighourhood\")\n\n# --- Edit 3: replace the \"Priority should...\" paragraph text and add a new\n#     paragraph after it with the expanded priority explanation ---\n$rng3 = $d.Content\n$rng3.Find.Text = \"Priority should be given to proximity, facilities, affordibility in that order\"\n$rng3.Find.Execute() | Out-Null\n$rng3.Text = \"From the already shorlisted neighbourhoods, I'll check which one has the least amount of rent.\"\n$rng3.InsertParagraphAfter()\n$rng3.Collapse($wdCollapseEnd)\n$rng3.MoveStart(1, 1) | Out-Null\n$rng3.InsertAfter(\"Priority should be given to proximity (will be taken care by foursquare radius variable), facilities (will be taken care by best fit neighourhoods), affordibility in that order (will be taken care by the rent data)\")\n"}
